$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4488.6665
$ws.Range("J18").Value = 6866.3335
$ws.Range("L18").Value = 6866.3335
$ws.Range("N18").Value = -7434.3335
$ws.Range("H125").Value = 2272.2727
$ws.Range("I125").Value = 3616
$ws.Range("K125").Value = 32544
$ws.Range("M125").Value = -30084
$ws.Range("H138").Value = 5684002.5
$ws.Range("J138").Value = 8930978
$ws.Range("L138").Value = 26792934
$ws.Range("N138").Value = -26803214

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9692.453
$ws.Range("I32").Value = 10264.809
$ws.Range("K32").Value = 10264.809
$ws.Range("M32").Value = -9977.808999999999
$ws.Range("H61").Value = 11630214
$ws.Range("I61").Value = 14287719
$ws.Range("K61").Value = 14287719
$ws.Range("M61").Value = -14287507
$ws.Range("H97").Value = 8617
$ws.Range("I97").Value = 11978.889
$ws.Range("J97").Value = 1052.75
$ws.Range("K97").Value = 11978.889
$ws.Range("L97").Value = 1052.75
$ws.Range("M97").Value = -11482.889
$ws.Range("N97").Value = -2044.75
$ws.Range("H136").Value = 11630214
$ws.Range("I136").Value = 14287719
$ws.Range("K136").Value = 42863157
$ws.Range("M136").Value = -42860607

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1812
$ws.Range("I20").Value = 2104
$ws.Range("K20").Value = 2104
$ws.Range("M20").Value = -1857
$ws.Range("H80").Value = 10938.895
$ws.Range("I80").Value = 25294.75
$ws.Range("J80").Value = 498.27274
$ws.Range("K80").Value = 25294.75
$ws.Range("L80").Value = 498.27274
$ws.Range("M80").Value = -24296.75
$ws.Range("N80").Value = -2494.27274
$ws.Range("H83").Value = 10938.895
$ws.Range("I83").Value = 25294.75
$ws.Range("J83").Value = 498.27274
$ws.Range("K83").Value = 126473.75
$ws.Range("L83").Value = 2491.3637
$ws.Range("M83").Value = -121481.75
$ws.Range("N83").Value = -12475.3637
$ws.Range("H86").Value = 16130887
$ws.Range("I86").Value = 1836.8422
$ws.Range("J86").Value = 41668548
$ws.Range("K86").Value = 1836.8422
$ws.Range("L86").Value = 41668548
$ws.Range("M86").Value = -713.8422
$ws.Range("N86").Value = -41670794
$ws.Range("H89").Value = 16130887
$ws.Range("I89").Value = 1836.8422
$ws.Range("J89").Value = 41668548
$ws.Range("K89").Value = 9184.210999999999
$ws.Range("L89").Value = 208342740
$ws.Range("M89").Value = -3568.210999999999
$ws.Range("N89").Value = -208353972
$ws.Range("H105").Value = 4194.972
$ws.Range("I105").Value = 2908.3845
$ws.Range("K105").Value = 2908.3845
$ws.Range("M105").Value = -1161.3845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 15513
$ws.Range("I41").Value = 2500
$ws.Range("J41").Value = 18766.25
$ws.Range("K41").Value = 2500
$ws.Range("L41").Value = 18766.25
$ws.Range("M41").Value = -2072
$ws.Range("N41").Value = -19622.25
$ws.Range("H50").Value = 11092
$ws.Range("J50").Value = 11092
$ws.Range("L50").Value = 11092
$ws.Range("N50").Value = -12342
$ws.Range("H51").Value = 993.3333
$ws.Range("I51").Value = 993.3333
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 993.3333
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -257.3333
$ws.Range("N51").ClearContents()
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H60").Value = 6676.5
$ws.Range("I60").Value = 2250
$ws.Range("J60").Value = 11103
$ws.Range("K60").Value = 2250
$ws.Range("L60").Value = 11103
$ws.Range("M60").Value = -1739
$ws.Range("N60").Value = -12125
$ws.Range("H61").Value = 993.3333
$ws.Range("I61").Value = 993.3333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 993.3333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -645.3333
$ws.Range("N61").ClearContents()
$ws.Range("H62").Value = 2381.9048
$ws.Range("I62").Value = 2372.9412
$ws.Range("J62").Value = 2420
$ws.Range("K62").Value = 2372.9412
$ws.Range("L62").Value = 2420
$ws.Range("M62").Value = -1748.9412
$ws.Range("N62").Value = -3668
$ws.Range("H65").Value = 2381.9048
$ws.Range("I65").Value = 2372.9412
$ws.Range("J65").Value = 2420
$ws.Range("K65").Value = 11864.706
$ws.Range("L65").Value = 12100
$ws.Range("M65").Value = -8744.706000000002
$ws.Range("N65").Value = -18340
$ws.Range("H74").Value = 33999.5
$ws.Range("J74").Value = 33999.5
$ws.Range("L74").Value = 33999.5
$ws.Range("N74").Value = -35747.5
$ws.Range("H77").Value = 33999.5
$ws.Range("J77").Value = 33999.5
$ws.Range("L77").Value = 101998.5
$ws.Range("N77").Value = -110734.5
$ws.Range("H105").Value = 2385
$ws.Range("I105").Value = 1318
$ws.Range("K105").Value = 1318
$ws.Range("M105").Value = 429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3465.0334
$ws.Range("I126").Value = 2061.9412
$ws.Range("J126").Value = 5299.846
$ws.Range("K126").Value = 6185.823600000001
$ws.Range("L126").Value = 15899.538
$ws.Range("M126").Value = -3715.823600000001
$ws.Range("N126").Value = -20839.538
$ws.Range("H132").Value = 4545.2583
$ws.Range("I132").Value = 3297.9443
$ws.Range("J132").Value = 6272.3076
$ws.Range("K132").Value = 9893.832900000001
$ws.Range("L132").Value = 18816.9228
$ws.Range("M132").Value = -7363.832900000001
$ws.Range("N132").Value = -23876.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4584.5713
$ws.Range("I7").Value = 4672.1875
$ws.Range("J7").Value = 4510.7896
$ws.Range("K7").Value = 4672.1875
$ws.Range("L7").Value = 4510.7896
$ws.Range("M7").Value = -4560.1875
$ws.Range("N7").Value = -4734.7896
$ws.Range("H22").Value = 1426.7273
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 1549.4
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 1549.4
$ws.Range("M22").Value = 95
$ws.Range("N22").Value = -2139.4
$ws.Range("H27").Value = 1426.7273
$ws.Range("I27").Value = 200
$ws.Range("J27").Value = 1549.4
$ws.Range("K27").Value = 200
$ws.Range("L27").Value = 1549.4
$ws.Range("M27").Value = -93
$ws.Range("N27").Value = -1763.4
$ws.Range("H46").Value = 1511.1111
$ws.Range("I46").Value = 1410
$ws.Range("J46").Value = 1637.5
$ws.Range("K46").Value = 1410
$ws.Range("L46").Value = 1637.5
$ws.Range("M46").Value = -1222
$ws.Range("N46").Value = -2013.5
$ws.Range("H82").Value = 1900.2174
$ws.Range("I82").Value = 1682.8125
$ws.Range("J82").Value = 2397.1428
$ws.Range("K82").Value = 1682.8125
$ws.Range("L82").Value = 2397.1428
$ws.Range("M82").Value = -1321.8125
$ws.Range("N82").Value = -3119.1428
$ws.Range("H85").Value = 1900.2174
$ws.Range("I85").Value = 1682.8125
$ws.Range("J85").Value = 2397.1428
$ws.Range("K85").Value = 1682.8125
$ws.Range("L85").Value = 2397.1428
$ws.Range("M85").Value = -434.8125
$ws.Range("N85").Value = -4893.1428
$ws.Range("H126").Value = 4584.5713
$ws.Range("I126").Value = 4672.1875
$ws.Range("J126").Value = 4510.7896
$ws.Range("K126").Value = 14016.5625
$ws.Range("L126").Value = 13532.3688
$ws.Range("M126").Value = -11546.5625
$ws.Range("N126").Value = -18472.3688
$ws.Range("H136").Value = 19238108
$ws.Range("I136").Value = 27780138
$ws.Range("K136").Value = 83340414
$ws.Range("M136").Value = -83337864

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 628
$ws.Range("I107").Value = 667.36
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 2002.08
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = -82.07999999999993
$ws.Range("N107").Value = -4740
$ws.Range("H113").Value = 1643.1578
$ws.Range("I113").Value = 563.375
$ws.Range("J113").Value = 2428.4546
$ws.Range("K113").Value = 1690.125
$ws.Range("L113").Value = 7285.3638
$ws.Range("M113").Value = 479.875
$ws.Range("N113").Value = -11625.3638
$ws.Range("H122").Value = 3023.75
$ws.Range("I122").Value = 2924.5833
$ws.Range("J122").Value = 3321.25
$ws.Range("K122").Value = 8773.749899999999
$ws.Range("L122").Value = 9963.75
$ws.Range("M122").Value = -6323.749899999999
$ws.Range("N122").Value = -14863.75
$ws.Range("H136").Value = 1069.9117
$ws.Range("I136").Value = 966.8570999999999
$ws.Range("K136").Value = 2900.5713
$ws.Range("M136").Value = -350.5712999999996
